$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this pushes the existing rows 16-46
# down to 17-47 (matching the new dimension A1:R47 from the diff).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the data added by the edit.
# (Columns A, B, C, E, F, G, H, I, R are constant across all data rows
# in this sheet.)
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44607
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = 100114002
$ws.Cells.Item(16, 7).Value = "Camote"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 40
$ws.Cells.Item(16, 11).Value = 18000
$ws.Cells.Item(16, 12).Value = 18000
$ws.Cells.Item(16, 13).Value = 18000
$ws.Cells.Item(16, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 900
$ws.Cells.Item(16, 17).Value = 20
$ws.Cells.Item(16, 18).Value = "Hortaliza"
